# Hortaliza, Femacal de La Calera - Acelga.xlsx
# Commit: "Fruta / hortaliza, semanal"
#
# A new weekly price-report row is inserted at row 153, pushing the
# existing rows 153-239 down to 154-240 (dimension grows from A1:R239
# to A1:R240). The new row carries a fresh sample for the same market
# (Femacal de La Calera / Coquimbo / Acelga / Primera).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 153:239 down to 154:240, opening up a blank row 153.
$ws.Rows("153:153").Insert()

# Fill in the newly inserted row 153 with the new observation.
$ws.Range("A153").Value = 3
$ws.Range("B153").Value = "Femacal de La Calera"
$ws.Range("C153").Value = "Coquimbo"
$ws.Range("D153").Value = 44529
$ws.Range("E153").Value = 5
$ws.Range("F153").Value = 100112009
$ws.Range("G153").Value = "Acelga"
$ws.Range("H153").Value = "Sin especificar"
$ws.Range("I153").Value = "Primera"
$ws.Range("J153").Value = 270
$ws.Range("K153").Value = 2000
$ws.Range("L153").Value = 2200
$ws.Range("M153").Value = 2104
$ws.Range("N153").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O153").Value = "Provincia de Quillota"
$ws.Range("P153").Value = 351
$ws.Range("Q153").Value = 6
$ws.Range("R153").Value = "Hortaliza"
